$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.710.47'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.28%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.387.96'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.48%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '503.65'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.91'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.51%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.550'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.73%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.389.77'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0971'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.57%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.149'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.64%  '
$ws.Range("E12").Value = '  -0.38%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.62'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.813.77'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.21%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '56.680.58'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.27%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.60'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.46%  '
$ws.Range("E17").Value = '  +0.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.422.72'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.59%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.19'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.05'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '308.49'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.26'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.53'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '67.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.98%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.52%  '
$ws.Range("E27").Value = '  -1.35%  '
$ws.Range("E28").Value = '  -0.97%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.40'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '175.67'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.68%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0721'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.73%  '
$ws.Range("E32").Value = '  -1.24%  '
$ws.Range("E33").Value = '  +0.62%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.86'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.73%  '
$ws.Range("E35").Value = '  +0.18%  '
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.81'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("E38").Value = '  -2.59%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.79'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.84%  '
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'OKB'
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.80'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.61%  '
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'SuiNetwork'
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.823'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.73%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.43'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '131.30'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.37'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.85'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.565'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '250.64'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0908'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.76%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0484'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0209'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.07'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.93%  '
